# Applies the "Finalizado artigo. Adicionado calculo das temperaturas de
# juncao." edit to Plan3: bumps Rjc (C21) from 0.25 to 0.5, fills in the
# measured Tc data for the "com 50CFM" block (rows 41-43) and builds four
# new Tc/Tj blocks (Hollowfin 50CFm / sem ventilacao / 25CFM / 10CFM) in
# rows 46-68, each with a junction-temperature (Tj) column computed from
# the Rjc*P offset.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Plan3")

# ---------------------------------------------------------------------
# Style helpers - mirror the 5 cell styles already used on this sheet so
# new cells look identical to their pre-existing neighbours (s2/s4/s11/
# s16/s17 in the original OOXML).
# ---------------------------------------------------------------------

function Set-BoxStyle($range) {
    # s="2" : General number format, thin box border, no special align.
    $range.NumberFormat = "General"
    $range.Borders.LineStyle = 1
}

function Set-DecimalStyle($range) {
    # s="4" : 0.00 number format, no border.
    $range.NumberFormat = "0.00"
    $range.Borders.LineStyle = -4142
}

function Set-BlankStyle($range) {
    # s="11": General, no border - used for empty filler cells.
    $range.ClearFormats()
}

function Set-SectionHeaderStyle($range) {
    # s="16": General, no border, horizontal-left align.
    $range.ClearFormats()
    $range.HorizontalAlignment = -4131
}

function Set-LabelStyle($range) {
    # s="17": General, thin box border, center/center/wrap align.
    $range.NumberFormat = "General"
    $range.Borders.LineStyle = 1
    $range.HorizontalAlignment = -4108
    $range.VerticalAlignment = -4108
    $range.WrapText = $true
}

# ---------------------------------------------------------------------
# 1) Rjc goes from 0.25 to 0.5 W/K - this alone ripples through the
#    existing Tc-max formulas in F20:H22 and the new Tj formulas below.
# ---------------------------------------------------------------------
$ws.Range("C21").Value = 0.5

# ---------------------------------------------------------------------
# 2) Row 40/41/42/43 block ("... com 50CFM") - add the measured Tc data
#    plus the new Tj (K:M) columns.
# ---------------------------------------------------------------------
Set-BoxStyle($ws.Range("K40:M40"))
$ws.Range("K40").Value = "17A"
$ws.Range("L40").Value = "24A"
$ws.Range("M40").Value = "31A"

$ws.Range("F41").Value = 71.2
$ws.Range("G41").Value = 81.8
$ws.Range("H41").Value = 90.7

$ws.Range("F42").Value = 81.8
$ws.Range("G42").Value = 95.9
$ws.Range("H42").Value = 110

$ws.Range("F43").Value = 92.4
$ws.Range("G43").Value = 110
$ws.Range("H43").Value = 129

Set-LabelStyle($ws.Range("J41:J43"))
$ws.Range("J41:J43").Merge()
$ws.Range("J41").Value = "Tj (ºC)"

Set-BoxStyle($ws.Range("K41:M43"))
$ws.Range("K41").Formula = '=F41+J15*$C$21'
$ws.Range("L41").Formula = '=G41+K15*$C$21'
$ws.Range("M41").Formula = '=H41+L15*$C$21'
$ws.Range("K42").Formula = '=F42+J16*$C$21'
$ws.Range("L42").Formula = '=G42+K16*$C$21'
$ws.Range("M42").Formula = '=H42+L16*$C$21'
$ws.Range("K43").Formula = '=F43+J17*$C$21'
$ws.Range("L43").Formula = '=G43+K17*$C$21'
$ws.Range("M43").Formula = '=H43+L17*$C$21'

# ---------------------------------------------------------------------
# Helper that builds one full "Tc (ºC) / Tj (ºC)" block:
#   row+0 : section title in column E (merged E:H not required - title
#           only lives in E, F:H just get the blank box style)
#   row+1 : "17A"/"24A"/"31A" headers in F:H and K:M
#   row+2..row+4 : Tc data (F:H) + Tj formulas (K:M), E/J labels merged
# ---------------------------------------------------------------------
function Build-Block($row, $titleStringIndex, $title, $tcVals) {
    $titleRow = $row
    $headerRow = $row + 1
    $r0 = $row + 2
    $r1 = $row + 3
    $r2 = $row + 4

    Set-SectionHeaderStyle($ws.Range("E" + $titleRow))
    $ws.Range("E" + $titleRow).Value = $title
    Set-BlankStyle($ws.Range("F" + $titleRow + ":H" + $titleRow))

    Set-BoxStyle($ws.Range("F" + $headerRow + ":H" + $headerRow))
    $ws.Range("F" + $headerRow).Value = "17A"
    $ws.Range("G" + $headerRow).Value = "24A"
    $ws.Range("H" + $headerRow).Value = "31A"
    Set-BoxStyle($ws.Range("K" + $headerRow + ":M" + $headerRow))
    $ws.Range("K" + $headerRow).Value = "17A"
    $ws.Range("L" + $headerRow).Value = "24A"
    $ws.Range("M" + $headerRow).Value = "31A"

    Set-LabelStyle($ws.Range("E" + $r0 + ":E" + $r2))
    $ws.Range("E" + $r0 + ":E" + $r2).Merge()
    $ws.Range("E" + $r0).Value = "Tc (ºC)"

    Set-LabelStyle($ws.Range("J" + $r0 + ":J" + $r2))
    $ws.Range("J" + $r0 + ":J" + $r2).Merge()
    $ws.Range("J" + $r0).Value = "Tj (ºC)"

    Set-BoxStyle($ws.Range("F" + $r0 + ":H" + $r2))
    Set-BoxStyle($ws.Range("K" + $r0 + ":M" + $r2))

    # Fill Tc measured values (may contain $null for intentionally blank cells).
    $rows = @($r0, $r1, $r2)
    for ($i = 0; $i -lt 3; $i++) {
        $rr = $rows[$i]
        if ($tcVals[$i][0] -ne $null) { $ws.Range("F" + $rr).Value = $tcVals[$i][0] }
        if ($tcVals[$i][1] -ne $null) { $ws.Range("G" + $rr).Value = $tcVals[$i][1] }
        if ($tcVals[$i][2] -ne $null) { $ws.Range("H" + $rr).Value = $tcVals[$i][2] }
    }

    # Tj formulas - each of the 3 data rows references the matching J/K/L
    # multiplier row (15, 16, 17) from the "Por IGBT (W)" table above.
    $mrows = @(15, 16, 17)
    for ($i = 0; $i -lt 3; $i++) {
        $rr = $rows[$i]
        $mr = $mrows[$i]
        $ws.Range("K" + $rr).Formula = '=F' + $rr + '+J' + $mr + '*$C$21'
        $ws.Range("L" + $rr).Formula = '=G' + $rr + '+K' + $mr + '*$C$21'
        $ws.Range("M" + $rr).Formula = '=H' + $rr + '+L' + $mr + '*$C$21'
    }
}

# ---------------------------------------------------------------------
# 3) Four new blocks: Hollowfin 50CFm / sem ventilacao / 25CFM / 10CFM.
# ---------------------------------------------------------------------
Build-Block 46 44 "Simulação Hollowfin com 50CFm" @(
    @(55.7, 58.5, 60.9),
    @(58.5, 62.3, 66.1),
    @(61.4, 66.1, 71.3)
)

Build-Block 52 45 "Simulação Hollowfin sem ventilação" @(
    @($null, 120, $null),
    @($null, $null, 162),
    @($null, $null, $null)
)

Build-Block 58 46 "Simulação Hollowfin com 25CFM" @(
    @($null, $null, $null),
    @($null, $null, $null),
    @($null, $null, 84.1)
)

Build-Block 64 47 "Simulação Hollowfin com 10CFM" @(
    @($null, $null, $null),
    @($null, $null, $null),
    @($null, $null, 116)
)

# ---------------------------------------------------------------------
# 4) Move the view: scroll down a bit and select H69 (one past the last
#    new row), matching the saved UI state of the edited workbook.
# ---------------------------------------------------------------------
$ws.Range("H69").Select()
